# "Generate Report for Handoff"
#
# The handback status for 27c80a73-16fb-4437-a628-5ab6f9ace938 is now stale
# (a newer commit exists upstream), so it flips back to "Ready for handoff"
# with a new error detail, while a654ead0-324b-4acb-baa0-a6aaf223040d keeps
# its "Handed back: in sync with en-US" status. Update the Overview sheet and
# both locale sheets (zh-cn, de-de) to reflect this, including the hyperlink
# display text and the widened "Error Detail" column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "a654ead0-324b-4acb-baa0-a6aaf223040d.md"
$ov.Range("B2").Value = "e2e\a654ead0-324b-4acb-baa0-a6aaf223040d.md"
$ov.Range("A3").Value = "27c80a73-16fb-4437-a628-5ab6f9ace938.md"
$ov.Range("B3").Value = "e2e\27c80a73-16fb-4437-a628-5ab6f9ace938.md"

$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-10-17 15:52:53"

foreach ($hl in $ov.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\a654ead0-324b-4acb-baa0-a6aaf223040d.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\27c80a73-16fb-4437-a628-5ab6f9ace938.md"
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "a654ead0-324b-4acb-baa0-a6aaf223040d.md"
$zh.Range("G2").Value = "a654ead0-324b-4acb-baa0-a6aaf223040d.c43ecf382617a9bfcc3ee3cc4c6c94a480119d18.zh-cn.xlf"
$zh.Range("I2").Value = "a654ead0-324b-4acb-baa0-a6aaf223040d.md"
$zh.Range("J2").Value = "a654ead0-324b-4acb-baa0-a6aaf223040d.c43ecf382617a9bfcc3ee3cc4c6c94a480119d18.zh-cn.xlf"

$zh.Range("A3").Value = "27c80a73-16fb-4437-a628-5ab6f9ace938.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "27c80a73-16fb-4437-a628-5ab6f9ace938.d65c5b990478e2d698490926e3efc92ea8caa28d.zh-cn.xlf"
$zh.Range("H3").Value = "2016-10-17 15:52:31"
$zh.Range("I3").Value = "27c80a73-16fb-4437-a628-5ab6f9ace938.md"
$zh.Range("J3").Value = "27c80a73-16fb-4437-a628-5ab6f9ace938.d65c5b990478e2d698490926e3efc92ea8caa28d.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c212fd0d7ec3c1bbcde43f01354f2d828af3889c/e2e/27c80a73-16fb-4437-a628-5ab6f9ace938.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8ea8804546ab7d78b908591a063b406c7e9ae699/e2e/27c80a73-16fb-4437-a628-5ab6f9ace938.md."

$zh.Columns.Item(16).ColumnWidth = 40

foreach ($hl in $zh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if (($addr -eq '$A$2') -or ($addr -eq '$I$2')) {
        $hl.TextToDisplay = "a654ead0-324b-4acb-baa0-a6aaf223040d.md"
    } elseif (($addr -eq '$A$3') -or ($addr -eq '$I$3')) {
        $hl.TextToDisplay = "27c80a73-16fb-4437-a628-5ab6f9ace938.md"
    }
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "a654ead0-324b-4acb-baa0-a6aaf223040d.md"
$de.Range("G2").Value = "a654ead0-324b-4acb-baa0-a6aaf223040d.c43ecf382617a9bfcc3ee3cc4c6c94a480119d18.de-de.xlf"
$de.Range("I2").Value = "a654ead0-324b-4acb-baa0-a6aaf223040d.md"
$de.Range("J2").Value = "a654ead0-324b-4acb-baa0-a6aaf223040d.c43ecf382617a9bfcc3ee3cc4c6c94a480119d18.de-de.xlf"

$de.Range("A3").Value = "27c80a73-16fb-4437-a628-5ab6f9ace938.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "27c80a73-16fb-4437-a628-5ab6f9ace938.d65c5b990478e2d698490926e3efc92ea8caa28d.de-de.xlf"
$de.Range("H3").Value = "2016-10-17 15:52:53"
$de.Range("I3").Value = "27c80a73-16fb-4437-a628-5ab6f9ace938.md"
$de.Range("J3").Value = "27c80a73-16fb-4437-a628-5ab6f9ace938.d65c5b990478e2d698490926e3efc92ea8caa28d.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c212fd0d7ec3c1bbcde43f01354f2d828af3889c/e2e/27c80a73-16fb-4437-a628-5ab6f9ace938.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8ea8804546ab7d78b908591a063b406c7e9ae699/e2e/27c80a73-16fb-4437-a628-5ab6f9ace938.md."

$de.Columns.Item(16).ColumnWidth = 40

foreach ($hl in $de.Hyperlinks) {
    $addr = $hl.Range.Address()
    if (($addr -eq '$A$2') -or ($addr -eq '$I$2')) {
        $hl.TextToDisplay = "a654ead0-324b-4acb-baa0-a6aaf223040d.md"
    } elseif (($addr -eq '$A$3') -or ($addr -eq '$I$3')) {
        $hl.TextToDisplay = "27c80a73-16fb-4437-a628-5ab6f9ace938.md"
    }
}
